$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> updated Notified Production (B) / Actual Production (C) values
# for rows whose production figures changed in addition to the date shift.
$bcMap = @{}
$bcMap[22] = @{B=17}
$bcMap[23] = @{B=20; C=5}
$bcMap[24] = @{B=24; C=20}
$bcMap[25] = @{B=29; C=35}
$bcMap[26] = @{B=140; C=63}
$bcMap[27] = @{B=149; C=92}
$bcMap[28] = @{B=162; C=128}
$bcMap[29] = @{B=175; C=156}
$bcMap[30] = @{B=378; C=220}
$bcMap[31] = @{B=426; C=275}
$bcMap[32] = @{B=446; C=319}
$bcMap[33] = @{B=466; C=357}
$bcMap[34] = @{B=647; C=377}
$bcMap[35] = @{B=664; C=395}
$bcMap[36] = @{B=686; C=455}
$bcMap[37] = @{B=708; C=501}
$bcMap[38] = @{B=836; C=565}
$bcMap[39] = @{B=852; C=609}
$bcMap[40] = @{B=871; C=648}
$bcMap[41] = @{B=884; C=652}
$bcMap[42] = @{B=992; C=670}
$bcMap[43] = @{B=1000; C=668}
$bcMap[44] = @{B=1008; C=692}
$bcMap[45] = @{B=1018; C=733}
$bcMap[46] = @{B=1090; C=759}
$bcMap[47] = @{B=1099; C=753}
$bcMap[48] = @{B=1104; C=703}
$bcMap[49] = @{B=1109; C=681}
$bcMap[50] = @{B=1104; C=684}
$bcMap[51] = @{B=1102; C=703}
$bcMap[52] = @{B=1096; C=699}
$bcMap[53] = @{B=1091; C=695}
$bcMap[54] = @{B=1042; C=678}
$bcMap[55] = @{B=1035; C=685}
$bcMap[56] = @{B=1026; C=717}
$bcMap[57] = @{B=1017; C=700}
$bcMap[58] = @{B=918; C=707}
$bcMap[59] = @{B=902; C=675}
$bcMap[60] = @{B=882; C=656}
$bcMap[61] = @{B=864; C=605}
$bcMap[62] = @{B=740; C=573}
$bcMap[63] = @{B=717; C=551}
$bcMap[64] = @{B=699; C=497}
$bcMap[65] = @{B=684; C=456}
$bcMap[66] = @{B=510; C=345}
$bcMap[67] = @{B=492; C=294}
$bcMap[68] = @{B=475; C=267}
$bcMap[69] = @{B=457; C=235}
$bcMap[70] = @{B=261; C=175}
$bcMap[71] = @{B=248; C=155}
$bcMap[72] = @{B=236; C=134}
$bcMap[73] = @{B=224; C=108}
$bcMap[74] = @{B=64; C=80}
$bcMap[75] = @{B=56; C=59}
$bcMap[76] = @{B=50; C=40}
$bcMap[77] = @{B=45; C=19}
$bcMap[78] = @{B=22; C=8}
$bcMap[79] = @{B=19; C=2}
$bcMap[80] = @{B=19; C=0}
$bcMap[81] = @{B=19; C=0}
$bcMap[82] = @{B=16; C=0}
$bcMap[83] = @{B=16}
$bcMap[84] = @{B=16}
$bcMap[85] = @{B=16}
$bcMap[86] = @{B=16}
$bcMap[118] = @{B=19}
$bcMap[119] = @{B=22; C=6}
$bcMap[120] = @{C=19}
$bcMap[121] = @{B=31; C=36}
$bcMap[122] = @{B=193; C=69}
$bcMap[123] = @{B=208; C=0}
$bcMap[124] = @{B=226; C=0}
$bcMap[125] = @{B=246; C=0}
$bcMap[126] = @{B=611; C=0}
$bcMap[127] = @{B=640; C=0}
$bcMap[128] = @{B=674; C=0}
$bcMap[129] = @{B=711; C=0}
$bcMap[130] = @{B=1041; C=0}
$bcMap[131] = @{B=1075; C=0}
$bcMap[132] = @{B=1111; C=0}
$bcMap[133] = @{B=1146; C=0}
$bcMap[134] = @{B=1368; C=0}
$bcMap[135] = @{B=1396}
$bcMap[136] = @{B=1424}
$bcMap[137] = @{B=1446}
$bcMap[138] = @{B=1589}
$bcMap[139] = @{B=1601}
$bcMap[140] = @{B=1612}
$bcMap[141] = @{B=1625}
$bcMap[142] = @{B=1703}
$bcMap[143] = @{B=1709}
$bcMap[144] = @{B=1720}
$bcMap[145] = @{B=1722}
$bcMap[146] = @{B=1766}
$bcMap[147] = @{B=1766}
$bcMap[148] = @{B=1762}
$bcMap[149] = @{B=1758}
$bcMap[150] = @{B=1706}
$bcMap[151] = @{B=1698}
$bcMap[152] = @{B=1687}
$bcMap[153] = @{B=1675}
$bcMap[154] = @{B=1561}
$bcMap[155] = @{B=1544}
$bcMap[156] = @{B=1528}
$bcMap[157] = @{B=1509}
$bcMap[158] = @{B=1335}
$bcMap[159] = @{B=1305}
$bcMap[160] = @{B=1279}
$bcMap[161] = @{B=1250}
$bcMap[162] = @{B=972}
$bcMap[163] = @{B=944}
$bcMap[164] = @{B=915}
$bcMap[165] = @{B=885}
$bcMap[166] = @{B=514}
$bcMap[167] = @{B=488}
$bcMap[168] = @{B=462}
$bcMap[169] = @{B=440}
$bcMap[170] = @{B=139}
$bcMap[171] = @{B=124}
$bcMap[172] = @{B=112}
$bcMap[173] = @{B=101}
$bcMap[174] = @{B=26}
$bcMap[175] = @{B=23}
$bcMap[177] = @{B=21}
$bcMap[178] = @{B=16}
$bcMap[179] = @{B=16}
$bcMap[180] = @{B=16}
$bcMap[181] = @{B=16}
$bcMap[182] = @{B=16}

for ($r = 2; $r -le 193; $r++) {
    # Shift the Timestamp (column A) forward by 10 days, preserving the
    # exact time-of-day fraction and the cell's existing date style.
    $oldDate = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r, 1).Value = $oldDate + 10

    # Update the Lookup text (column E) to reflect the new date while
    # keeping the same quarter-hour index suffix.
    $oldLookup = $ws.Cells.Item($r, 5).Value()
    $newLookup = $oldLookup.Replace("28.04.2025", "08.05.2025").Replace("29.04.2025", "09.05.2025")
    $ws.Cells.Item($r, 5).Value = $newLookup

    # Apply any production-figure corrections for this row, if present.
    if ($bcMap.ContainsKey($r)) {
        $entry = $bcMap[$r]
        if ($entry.ContainsKey("B")) {
            $ws.Cells.Item($r, 2).Value = $entry.B
        }
        if ($entry.ContainsKey("C")) {
            $ws.Cells.Item($r, 3).Value = $entry.C
        }
    }
}
